$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1446896", "Asif Newaz", "2025-01-26", "14:46:01"),
    @("1446896", "Asif Newaz", "2025-01-26", "14:55:47"),
    @("1446896", "Asif Newaz", "2025-01-26", "14:56:55"),
    @("1446896", "Asif Newaz", "2025-01-26", "15:01:01"),
    @("1446896", "Asif Newaz", "2025-01-26", "15:07:06"),
    @("1420867", "Mohsina Binte Asad", "2025-01-26", "15:10:34"),
    @("1446896", "Asif Newaz", "2025-01-26", "15:10:48"),
    @("1420867", "Mohsina Binte Asad", "2025-01-26", "15:13:07"),
    @("1446896", "Asif Newaz", "2025-01-26", "15:13:38"),
    @("121212", "Mukit", "2025-01-26", "15:15:27"),
    @("1446896", "Asif Newaz", "2025-01-26", "15:15:53"),
    @("13131313", "Kamrul Hasan", "2025-01-26", "15:18:21"),
    @("13131313", "Kamrul Hasan", "2025-01-26", "15:19:05"),
    @("13131313", "Kamrul Hasan", "2025-01-26", "15:20:47"),
    @("13131313", "Kamrul Hasan", "2025-01-26", "15:27:12"),
    @("13131313", "Kamrul Hasan", "2025-01-26", "15:29:03"),
    @("1446896", "Asif Newaz", "2025-01-26", "15:32:58"),
    @("1446896", "Asif Newaz", "2025-01-26", "15:36:47"),
    @("13131313", "Kamrul Hasan", "2025-01-26", "16:45:04")
)

$row = 6
foreach ($entry in $data) {
    $ws.Range("A${row}:D${row}").NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
